$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 160: date shifts 44383 -> 44509, other columns unchanged ---
$ws.Range("D160").Value = 44509

# --- Row 161: date shifts 44383 -> 44509, other columns unchanged ---
$ws.Range("D161").Value = 44509

# --- Row 162: date 44273 -> 44383, prices/origin change ---
$ws.Range("D162").Value = 44383
$ws.Range("J162").Value = 1000
$ws.Range("K162").Value = 700
$ws.Range("L162").Value = 800
$ws.Range("M162").Value = 750
$ws.Range("O162").Value = "Región Metropolitana"
$ws.Range("P162").Value = 750

# --- Row 163: date 44273 -> 44383, prices/origin change ---
$ws.Range("D163").Value = 44383
$ws.Range("J163").Value = 500
$ws.Range("K163").Value = 600
$ws.Range("L163").Value = 600
$ws.Range("M163").Value = 600
$ws.Range("O163").Value = "Región Metropolitana"
$ws.Range("P163").Value = 600

# --- Row 164: date 44491 -> 44273, J and origin change ---
$ws.Range("D164").Value = 44273
$ws.Range("J164").Value = 800
$ws.Range("O164").Value = "Región del Maule"

# --- Row 165: date 44491 -> 44273, J and origin change ---
$ws.Range("D165").Value = 44273
$ws.Range("J165").Value = 400
$ws.Range("O165").Value = "Región del Maule"

# --- New row 166 ---
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = "Vega Monumental Concepción"
$ws.Range("C166").Value = "Bíobío"
$ws.Range("D166").Value = 44491
$ws.Range("D166").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = 100112008
$ws.Range("G166").Value = "Coliflor"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 1000
$ws.Range("K166").Value = 800
$ws.Range("L166").Value = 900
$ws.Range("M166").Value = 850
$ws.Range("N166").Value = "$/unidad"
$ws.Range("O166").Value = "Región Metropolitana"
$ws.Range("P166").Value = 850
$ws.Range("Q166").Value = 1
$ws.Range("R166").Value = "Hortaliza"

# --- New row 167 ---
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = "Vega Monumental Concepción"
$ws.Range("C167").Value = "Bíobío"
$ws.Range("D167").Value = 44491
$ws.Range("D167").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = 100112008
$ws.Range("G167").Value = "Coliflor"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Segunda"
$ws.Range("J167").Value = 500
$ws.Range("K167").Value = 700
$ws.Range("L167").Value = 700
$ws.Range("M167").Value = 700
$ws.Range("N167").Value = "$/unidad"
$ws.Range("O167").Value = "Región Metropolitana"
$ws.Range("P167").Value = 700
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
